# "Orthopedics United" roster sheet update.
#
# The fantasy roster table (A:Oyuncu Adı / B:Pozisyon / C:Takım) is updated:
#   - the row for "Khris Middleton" (previously the very last row, 19,
#     SF / Milwaukee Bucks) is moved up to sit right after "Buddy Hield",
#   - "Pascal Siakam" moves down to sit after "Jalen Green",
#   - "Aaron Wiggins" is removed entirely,
#   - "Chet Holmgren" / "Jalen Suggs" / "Deni Avdija" are re-ordered,
#   - the now-empty trailing row (old row 19) disappears, shrinking the
#     table from 18 data rows to 17 (A1:C19 -> A1:C18).
#
# Simplest reliable way to reproduce this with COM automation: delete the
# last row outright (the table only shrinks by one row), then rewrite the
# data rows 2-18 with the corrected roster content, in final order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing row - the table goes from 19 rows (incl. header) to 18.
$ws.Rows(19).Delete()

# Re-populate rows 2-18 with the corrected roster, in the new order.
$roster = @(
    @("Russell Westbrook",   "PG,SG",   "Denver Nuggets"),
    @("Dejounte Murray",     "PG,SG",   "New Orleans Pelicans"),
    @("Keon Johnson",        "PG,SG",   "Brooklyn Nets"),
    @("Chris Paul",          "PG",      "San Antonio Spurs"),
    @("Jaylen Brown",        "SG,SF",   "Boston Celtics"),
    @("Paolo Banchero",      "SF,PF",   "Orlando Magic"),
    @("Buddy Hield",         "SG,SF",   "Golden State Warriors"),
    @("Khris Middleton",     "SF",      "Milwaukee Bucks"),
    @("Nikola Jokic",        "C",       "Denver Nuggets"),
    @("Rudy Gobert",         "C",       "Minnesota Timberwolves"),
    @("Jakob Poeltl",        "C",       "Toronto Raptors"),
    @("Jalen Green",         "PG,SG",   "Houston Rockets"),
    @("Pascal Siakam",       "SF,PF,C", "Indiana Pacers"),
    @("Bennedict Mathurin",  "SG,SF",   "Indiana Pacers"),
    @("Chet Holmgren",       "PF,C",    "Oklahoma City Thunder"),
    @("Jalen Suggs",         "PG,SG",   "Orlando Magic"),
    @("Deni Avdija",         "SF,PF",   "Portland Trail Blazers")
)

$row = 2
foreach ($player in $roster) {
    $ws.Range("A$row").Value = $player[0]
    $ws.Range("B$row").Value = $player[1]
    $ws.Range("C$row").Value = $player[2]
    $row++
}
